$d = $word.ActiveDocument

# 1. Update the letter date: "September 19, 2025" -> "September 21, 2025"
$d.Content.Find.Execute("September 19, 2025", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "September 21, 2025", 2) | Out-Null

# 2. Split the mailing address paragraph ("969 Story Road, San Jose CA 95122")
#    into two paragraphs: "969 Story Road" and a new "San Jose, CA 95122" line.
#    Only the mailing-address occurrence (not the "PROPERTY ADDRESS:" table entry)
#    is touched, so locate it via the Paragraphs collection rather than a
#    document-wide Find/Replace.
$addrPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.TrimEnd([char]13, [char]7) -eq "969 Story Road, San Jose CA 95122") {
        $addrPara = $para
        break
    }
}

if ($addrPara -ne $null) {
    $full = $addrPara.Range
    # Trim off the trailing paragraph mark so only the visible text is replaced.
    $textRange = $d.Range($full.Start, $full.End - 1)
    $textRange.Text = "969 Story Road"

    # Re-fetch the (now shorter) paragraph and append a new paragraph after it,
    # inheriting the same paragraph/run formatting automatically.
    $addrPara = $d.Paragraphs.Item($i)
    $newPara = $addrPara.Range.InsertParagraphAfter()

    $cityPara = $d.Paragraphs.Item($i + 1)
    $cityPara.Range.Text = "San Jose, CA 95122"
}

# 3. Remove the extra empty "No Spacing" paragraph that sits right after the
#    "...Board of Directors" signature line.
for ($j = 1; $j -le $d.Paragraphs.Count; $j++) {
    $para = $d.Paragraphs.Item($j)
    if ($para.Range.Text.TrimEnd([char]13, [char]7) -eq "Vietnam Town Condominium Owners Association Board of Directors") {
        $nextPara = $d.Paragraphs.Item($j + 1)
        if ($nextPara.Range.Text.TrimEnd([char]13, [char]7) -eq "") {
            $nextPara.Range.Delete()
        }
        break
    }
}
